# feat(get merch): naming convension
#
# Appends a new order (rows 24-34 on "orders", row 4 on "customer_info")
# for the same customer (Brock Tomlinson), identified by a new Order Id
# guid "ea52f74b-7b77-4248-9773-3898b445486d".

$wb = $excel.ActiveWorkbook

$ordersSheet = $wb.Worksheets.Item("orders")
$customerSheet = $wb.Worksheets.Item("customer_info")

$orderId = "ea52f74b-7b77-4248-9773-3898b445486d"

# --- "orders" sheet: 11 new line items for the new order, rows 24-34 ---
$items = @(
    @("Floof Hoodie", "L", 1, "Black", 49.9900016784668),
    @("Floof CrewNeck", "XL", 1, "Black", 39.9900016784668),
    @("Floof CrewNeck by Brock the One and Only Rockstar", "XL", 1, "Black", 39.9900016784668),
    @("Floof CrewNeck", "XL", 1, "Black", 39.9900016784668),
    @("Floof CrewNeck", "XL", 1, "Black", 39.9900016784668),
    @("Floof CrewNeck", "XL", 1, "Black", 39.9900016784668),
    @("Floof CrewNeck", "XL", 1, "Black", 39.9900016784668),
    @("Floof CrewNeck", "XL", 1, "Black", 39.9900016784668),
    @("Floof CrewNeck", "XL", 1, "Black", 39.9900016784668),
    @("Floof CrewNeck", "XL", 1, "Black", 39.9900016784668),
    @("Floof CrewNeck", "XL", 1, "Black", 39.9900016784668)
)

$row = 24
foreach ($item in $items) {
    $ordersSheet.Cells.Item($row, 1).Value = $orderId
    $ordersSheet.Cells.Item($row, 2).Value = $item[0]
    $ordersSheet.Cells.Item($row, 3).Value = $item[1]
    $ordersSheet.Cells.Item($row, 4).Value = $item[2]
    $ordersSheet.Cells.Item($row, 5).Value = $item[3]
    $ordersSheet.Cells.Item($row, 6).Value = $item[4]
    $row = $row + 1
}

# --- "customer_info" sheet: one new row (row 4) duplicating Brock's info ---
$customerSheet.Cells.Item(4, 1).Value = $orderId
$customerSheet.Cells.Item(4, 2).Value = "brock.tomlinson@ucalgarybaja.ca"
$customerSheet.Cells.Item(4, 3).Value = "2509466196"
$customerSheet.Cells.Item(4, 4).Value = "Brock Tomlinson"
$customerSheet.Cells.Item(4, 5).Value = "Software"
$customerSheet.Cells.Item(4, 6).Value = 472.38446044921875
$customerSheet.Cells.Item(4, 7).Value = "Yes"
$customerSheet.Cells.Item(4, 9).Value = ""
$customerSheet.Cells.Item(4, 10).Value = ""
$customerSheet.Cells.Item(4, 11).Value = ""
$customerSheet.Cells.Item(4, 12).Value = ""
$customerSheet.Cells.Item(4, 13).Value = ""
$customerSheet.Cells.Item(4, 14).Value = "Canada"
$customerSheet.Cells.Item(4, 15).Value = ""
$customerSheet.Cells.Item(4, 16).Value = ""
$customerSheet.Cells.Item(4, 17).Value = "Notes"
